# Update "想去人数" (number of people interested) counts for a few events.
# Sheet "展览" (sheet1): F3 217->219, F4 836->840, F6 28->29
# Sheet "全部类型" (sheet4): F4 217->219, F5 836->840, F7 28->29

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 219
$wsExhibit.Range("F4").Value = 840
$wsExhibit.Range("F6").Value = 29

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 219
$wsAll.Range("F5").Value = 840
$wsAll.Range("F7").Value = 29
